# Updates cryptos list values per the scraped diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "54.166.84"
$ws.Range("E2").Value = "  -3.66%  "

# Row 3
$ws.Range("D3").Value = "2.269.50"
$ws.Range("E3").Value = "  -4.12%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'490.90"
$ws.Range("E5").Value = "  -2.92%  "

# Row 6
$ws.Range("D6").Value = "'126.87"
$ws.Range("E6").Value = "  -2.45%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.63%  "

# Row 8
$ws.Range("E8").Value = "  -2.86%  "

# Row 9
$ws.Range("D9").Value = "2.268.45"
$ws.Range("E9").Value = "  -4.43%  "

# Row 10
$ws.Range("D10").Value = "'0.0932"
$ws.Range("E10").Value = "  -5.43%  "

# Row 12
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  -0.82%  "

# Row 13
$ws.Range("D13").Value = "'4.66"
$ws.Range("E13").Value = "  -4.52%  "

# Row 14
$ws.Range("D14").Value = "2.666.47"
$ws.Range("E14").Value = "  -4.27%  "

# Row 15
$ws.Range("D15").Value = "'21.45"
$ws.Range("E15").Value = "  -1.47%  "

# Row 16
$ws.Range("D16").Value = "54.079.56"
$ws.Range("E16").Value = "  -3.80%  "

# Row 17
$ws.Range("E17").Value = "  -3.46%  "

# Row 18
$ws.Range("D18").Value = "2.250.75"
$ws.Range("E18").Value = "  -5.61%  "

# Row 19
$ws.Range("D19").Value = "'9.78"
$ws.Range("E19").Value = "  -2.19%  "

# Row 20
$ws.Range("D20").Value = "'4.03"
$ws.Range("E20").Value = "  -0.02%  "

# Row 21
$ws.Range("D21").Value = "'297.55"
$ws.Range("E21").Value = "  -3.90%  "

# Row 22
$ws.Range("D22").Value = "'6.24"
$ws.Range("E22").Value = "  -0.40%  "

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("D24").Value = "'63.64"
$ws.Range("E24").Value = "  -2.95%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.32%  "

# Row 26
$ws.Range("D26").Value = "'0.373"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("E27").Value = "  -0.44%  "

# Row 28
$ws.Range("D28").Value = "2.332.82"
$ws.Range("E28").Value = "  -5.98%  "

# Row 29
$ws.Range("D29").Value = "'7.10"
$ws.Range("E29").Value = "  -1.18%  "

# Row 30
$ws.Range("D30").Value = "'162.92"
$ws.Range("E30").Value = "  -5.85%  "

# Row 31
$ws.Range("E31").Value = "  -2.96%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0678"
$ws.Range("E32").Value = "  -4.68%  "

# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.02%  "

# Row 34
$ws.Range("D34").Value = "'5.79"
$ws.Range("E34").Value = "  -0.96%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.22%  "

# Row 36
$ws.Range("E36").Value = "  -0.47%  "

# Row 37
$ws.Range("E37").Value = "  -0.72%  "

# Row 38
$ws.Range("E38").Value = "  -0.05%  "

# Row 39
$ws.Range("D39").Value = "'0.836"
$ws.Range("E39").Value = "  +1.44%  "

# Row 40
$ws.Range("D40").Value = "'3.62"
$ws.Range("E40").Value = "  -1.52%  "

# Row 41
$ws.Range("D41").Value = "'35.44"
$ws.Range("E41").Value = "  -2.25%  "

# Row 42
$ws.Range("E42").Value = "  +0.86%  "

# Row 43
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("E44").Value = "  -1.35%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'126.00"
$ws.Range("E45").Value = "  +0.51%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.81"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47
$ws.Range("E47").Value = "  -1.03%  "

# Row 48
$ws.Range("D48").Value = "'242.19"
$ws.Range("E48").Value = "  +2.22%  "

# Row 49
$ws.Range("E49").Value = "  -3.23%  "

# Row 50
$ws.Range("E50").Value = "  -0.83%  "

# Row 51
$ws.Range("E51").Value = "  -1.70%  "

# Reset style on cells that were numeric-looking text, so Excel
# does not persist an extra quote-prefix style on them.
foreach ($addr in @("D5","D6","D7","D10","D12","D13","D15","D19","D20","D21","D22","D23","D24","D25","D26","D29","D30","D33","D34","D35","D39","D40","D41","D45","D46","D48")) {
    $ws.Range($addr).Style = "Normal"
}
